$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New value added in previously-empty cell C2
$ws.Range("C2").Value = -0.1474887606676356

# Updated values in column E (rows 2-19)
$ws.Range("E2").Value  = 1.265019766896436
$ws.Range("E3").Value  = 1.858657482882586
$ws.Range("E4").Value  = -3.246097549514837
$ws.Range("E5").Value  = 1.609625625599986
$ws.Range("E6").Value  = 1.609625625600009
$ws.Range("E7").Value  = 3.238605209599998
$ws.Range("E8").Value  = 2.777885851461503
$ws.Range("E9").Value  = 2.436566844071941
$ws.Range("E10").Value = 1.694971351092267
$ws.Range("E11").Value = 1.216098605743365
$ws.Range("E12").Value = 1.784618024189011
$ws.Range("E13").Value = 4.887093273600018
$ws.Range("E14").Value = -12.1986023424
$ws.Range("E15").Value = 12.1815000816919
$ws.Range("E16").Value = 5.870037016039187
$ws.Range("E17").Value = -1.24582517146522
$ws.Range("E18").Value = -0.1040473946152809
$ws.Range("E19").Value = 0.8749339604052775

# Updated values in column C (rows 3-19)
$ws.Range("C3").Value  = 0.337821977117625
$ws.Range("C4").Value  = -1.890773121057054
$ws.Range("C5").Value  = 1.89156560755015
$ws.Range("C6").Value  = 1.113165545862094
$ws.Range("C7").Value  = 1.070385798714391
$ws.Range("C8").Value  = 1.384186838979806
$ws.Range("C9").Value  = 2.349355943833098
$ws.Range("C10").Value = 1.786425635558397
$ws.Range("C11").Value = 1.331333081915509
$ws.Range("C12").Value = 1.282262557986447
$ws.Range("C13").Value = 2.247109253368285
$ws.Range("C14").Value = -4.247034401476801
$ws.Range("C15").Value = -2.608215948579529
$ws.Range("C16").Value = 4.863085601670813
$ws.Range("C17").Value = -1.44371442952016
$ws.Range("C18").Value = 0.06625622369935691
$ws.Range("C19").Value = 0.9919038146506631
